$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cw2")

# --- Row 4: L1 measurements ---
$ws.Range("D4").Value = 42.5
$ws.Range("E4").Value = 11
$ws.Range("F4").Value = 38
$ws.Range("G4").Value = 32
$ws.Range("H4").Value = 42

$ws.Range("M4").Value = 124
$ws.Range("N4").Value = 81
$ws.Range("O4").Value = 75
$ws.Range("P4").Value = 66
$ws.Range("Q4").Value = 87

# --- Row 5: L2 measurements ---
$ws.Range("D5").Value = 162
$ws.Range("E5").Value = 140
$ws.Range("F5").Value = 158
$ws.Range("G5").Value = 151.5
$ws.Range("H5").Value = 161

$ws.Range("M5").Value = 86
$ws.Range("N5").Value = 136
$ws.Range("O5").Value = 104
$ws.Range("P5").Value = 72
$ws.Range("Q5").Value = 138

# --- Row 6: phi_i formulas (N6:Q6 direction flipped vs before) ---
$ws.Range("N6").Formula = "=N5-N4"
$ws.Range("O6").Formula = "=O5-O4"
$ws.Range("P6").Formula = "=P5-P4"
$ws.Range("Q6").Formula = "=Q5-Q4"

# --- Row 7: B measurements ---
$ws.Range("M7").Value = 222
$ws.Range("N7").Value = 103
$ws.Range("O7").Value = 189
$ws.Range("P7").Value = 55
$ws.Range("Q7").Value = 48.5

# --- Row 8: beta2 measurements ---
$ws.Range("M8").Value = 245
$ws.Range("N8").Value = 6.5
$ws.Range("O8").Value = 106
$ws.Range("P8").Value = 134.5
$ws.Range("Q8").Value = 138

# --- Row 9: beta formulas (direction/order changed) ---
$ws.Range("M9").Formula = "=M8-M7"
$ws.Range("N9").Formula = "=N7-N8"
$ws.Range("O9").Formula = "=O7-O8"
$ws.Range("P9").Formula = "=P8-P7"
$ws.Range("Q9").Formula = "=Q8-Q7"

# --- Row 12: new derived cells D12, F12, N12, P12 ---
$ws.Range("D12").Formula = "=I12"
$ws.Range("F12").Formula = "=I12*57.2957795"
$ws.Range("N12").Formula = "=S12"
$ws.Range("P12").Formula = "=S12*57.2957795"

# --- Row 13: new derived cells D13, F13 ---
$ws.Range("D13").Formula = "=I13"
$ws.Range("F13").Formula = "=I13*57.2957795"

# --- Row 14: new derived cells D14, F14 ---
$ws.Range("D14").Formula = "=I14"
$ws.Range("F14").Formula = "=I14*57.2957795"

# --- Row 15: new derived cells D15, F15 ---
$ws.Range("D15").Formula = "=I15"
$ws.Range("F15").Formula = "=I15*57.2957795"

# --- Row 24: B24 / D24 labels "n=" and formulas for C24 / E24 restated ---
$ws.Range("B24").Value = "n="
$ws.Range("D24").Value = "n="
$ws.Range("C24").Formula = "=(SIN(0.5*(N12+D14)))/(SIN(0.5*D14))"
$ws.Range("E24").Formula = "=(SIN(0.5*(P12+F14)))/(SIN(0.5*F14))"

# --- Remove the old "n=" labels from C28 and C29 (moved to B24/D24) ---
$ws.Range("C28").ClearContents()
$ws.Range("C29").ClearContents()

# --- Row 34: P34 formula stays the same text but now resolves (no /0 error) ---
$ws.Range("P34").Formula = "=(ABS((SIN(P12/2))/(2*(SIN(F14/2)^2))))*(ABS(F15))+(ABS((COS((F14+P12)/2))/(2*SIN(F14/2))))*(ABS(S28))"

$wb.Save()
